$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$k2 = @'


<p>Welcome to your first data-based quiz on data visualization! For data-based quizzes, you will need to import data into R as you follow the questions.</p>
<p>Note that each time you retry a quiz, the dataset is slightly modified.</p>
<div id="hiv-dataset" class="section level1">
<h1>HIV dataset</h1>
<p>Here, you will analyze parts of a dataset documenting global HIV prevalence from 1990 to 2011.</p>
<p>Click <a href="https://drive.google.com/file/d/1DjcpRAWXaWZECRpUZIKYR01c7wh7EBlG/view?usp=drivesdk" target="_blank"><strong>here</strong></a> to view and download the data. Or import it directly into R with the code below:</p>
<pre><code>if(!require(pacman)) install.packages(&quot;pacman&quot;)   
pacman::p_load(rio)   
hiv_prevalence &lt;- import(&quot;https://docs.google.com/uc?id=1DjcpRAWXaWZECRpUZIKYR01c7wh7EBlG&amp;export=download&quot;,
 format = &quot;csv&quot;,
 setclass = &quot;tibble&quot;)  </code></pre>
<p>Here are the top 6 rows of <code>hiv_prevalence</code> after import:</p>
<pre class="r"><code>head(hiv_prevalence)</code></pre>
<pre><code>## # A tibble: 6 × 4
##   country year  total_cases population
##   &lt;fct&gt;   &lt;fct&gt;       &lt;dbl&gt;      &lt;dbl&gt;
## 1 Brazil  1995       350000  161890816
## 2 Brazil  1996       370000  164392423
## 3 Brazil  1997       390000  166925457
## 4 Brazil  1998       410000  169472347
## 5 Brazil  1999       430000  172006362
## 6 Brazil  2000       440000  174504898</code></pre>
<pre class="r"><code>summary(hiv_prevalence)</code></pre>
<pre><code>##      country        year     total_cases     
##  Brazil  :15   1995   : 3   Min.   : 230000  
##  Cameroon:15   1996   : 3   1st Qu.: 440000  
##  India   :15   1997   : 3   Median : 490000  
##                1998   : 3   Mean   :1053111  
##                1999   : 3   3rd Qu.:2000000  
##                2000   : 3   Max.   :2600000  
##                (Other):27                    
##    population       
##  Min.   :1.393e+07  
##  1st Qu.:1.861e+07  
##  Median :1.794e+08  
##  Mean   :4.236e+08  
##  3rd Qu.:1.008e+09  
##  Max.   :1.190e+09  
## </code></pre>
<p><strong>Complete the code below to create a scatter plot of HIV incidence over time.</strong></p>
<p>Hint: You plot should look like this:</p>
<p><img src="https://drive.google.com/uc?export=view&id=1SAZ98XmvjVMufT4pmBI-NQS36Ej0iL-L" width="480" /></p>
<pre class="r"><code>ggplot(data = {hiv_prevalence}, 
       {mapping} = aes({x} = year,
                       y = {total_cases})) {+}
  geom{_point}()</code></pre>
<hr />
<p>Make this plot and save it as an R object:
<img src="https://drive.google.com/uc?export=view&id=1xc85TX30oiQBEmia9Evou-tG8b9NBbPh" width="672" /></p>
<p>Put your plot object through the <code>ggplot_digest()</code> function and enter the resulting character string here: {9fb1cf69b4fde79134ffb871ae01bf56}</p>
</div>
'@

$k3 = @'


<p>Welcome to your first data-based quiz on data visualization! For data-based quizzes, you will need to import data into R as you follow the questions.</p>
<p>Note that each time you retry a quiz, the dataset is slightly modified.</p>
<div id="hiv-dataset" class="section level1">
<h1>HIV dataset</h1>
<p>Here, you will analyze parts of a dataset documenting global HIV prevalence from 1990 to 2011.</p>
<p>Click <a href="https://drive.google.com/file/d/10UHloNs7wu1ES2aTJHsiC-iOBP615vUu/view?usp=drivesdk" target="_blank"><strong>here</strong></a> to view and download the data. Or import it directly into R with the code below:</p>
<pre><code>if(!require(pacman)) install.packages(&quot;pacman&quot;)   
pacman::p_load(rio)   
hiv_prevalence &lt;- import(&quot;https://docs.google.com/uc?id=10UHloNs7wu1ES2aTJHsiC-iOBP615vUu&amp;export=download&quot;,
 format = &quot;csv&quot;,
 setclass = &quot;tibble&quot;)  </code></pre>
<p>Here are the top 6 rows of <code>hiv_prevalence</code> after import:</p>
<pre class="r"><code>head(hiv_prevalence)</code></pre>
<pre><code>## # A tibble: 6 × 4
##   country year  total_cases population
##   &lt;fct&gt;   &lt;fct&gt;       &lt;dbl&gt;      &lt;dbl&gt;
## 1 Ghana   1995       170000   16760926
## 2 Ghana   1996       190000   17169151
## 3 Ghana   1997       200000   17568461
## 4 Ghana   1998       220000   17968830
## 5 Ghana   1999       230000   18384302
## 6 Ghana   2000       240000   18825034</code></pre>
<pre class="r"><code>summary(hiv_prevalence)</code></pre>
<pre><code>##      country        year     total_cases    
##  Ghana   :15   1995   : 3   Min.   :170000  
##  Thailand:15   1996   : 3   1st Qu.:250000  
##  Zambia  :15   1997   : 3   Median :610000  
##                1998   : 3   Mean   :570000  
##                1999   : 3   3rd Qu.:790000  
##                2000   : 3   Max.   :950000  
##                (Other):27                   
##    population      
##  Min.   : 8841338  
##  1st Qu.:11781612  
##  Median :19786307  
##  Mean   :31334510  
##  3rd Qu.:60903042  
##  Max.   :66277335  
## </code></pre>
<p><strong>Complete the code below to create a scatter plot of HIV incidence over time.</strong></p>
<p>Hint: You plot should look like this:</p>
<p><img src="https://drive.google.com/uc?export=view&id=1MJIjGZcM15tJv51ogPh3MlbYxFnpmKip" width="480" /></p>
<pre class="r"><code>ggplot(data = {hiv_prevalence}, 
       {mapping} = aes({x} = year,
                       y = {total_cases})) {+}
  geom{_point}()</code></pre>
<hr />
<p>Make this plot and save it as an R object:
<img src="https://drive.google.com/uc?export=view&id=1JsmOeM4rAwHylXTniph0axZOmjUHFr9y" width="672" /></p>
<p>Put your plot object through the <code>ggplot_digest()</code> function and enter the resulting character string here: {865fce3a2bd569ad8ee8629e845e9b1e}</p>
</div>
'@

$k4 = @'


<p>Welcome to your first data-based quiz on data visualization! For data-based quizzes, you will need to import data into R as you follow the questions.</p>
<p>Note that each time you retry a quiz, the dataset is slightly modified.</p>
<div id="hiv-dataset" class="section level1">
<h1>HIV dataset</h1>
<p>Here, you will analyze parts of a dataset documenting global HIV prevalence from 1990 to 2011.</p>
<p>Click <a href="https://drive.google.com/file/d/1rIYcr4AuMCUmcQk8ED07tWhHx_WG0uNx/view?usp=drivesdk" target="_blank"><strong>here</strong></a> to view and download the data. Or import it directly into R with the code below:</p>
<pre><code>if(!require(pacman)) install.packages(&quot;pacman&quot;)   
pacman::p_load(rio)   
hiv_prevalence &lt;- import(&quot;https://docs.google.com/uc?id=1rIYcr4AuMCUmcQk8ED07tWhHx_WG0uNx&amp;export=download&quot;,
 format = &quot;csv&quot;,
 setclass = &quot;tibble&quot;)  </code></pre>
<p>Here are the top 6 rows of <code>hiv_prevalence</code> after import:</p>
<pre class="r"><code>head(hiv_prevalence)</code></pre>
<pre><code>## # A tibble: 6 × 4
##   country  year  total_cases population
##   &lt;fct&gt;    &lt;fct&gt;       &lt;dbl&gt;      &lt;dbl&gt;
## 1 Ethiopia 1995       910000   57023519
## 2 Ethiopia 1996      1000000   58815116
## 3 Ethiopia 1997      1100000   60584273
## 4 Ethiopia 1998      1200000   62353942
## 5 Ethiopia 1999      1200000   64158887
## 6 Ethiopia 2000      1300000   66024199</code></pre>
<pre class="r"><code>summary(hiv_prevalence)</code></pre>
<pre><code>##      country        year     total_cases     
##  Ethiopia:15   1995   : 3   Min.   : 610000  
##  India   :15   1996   : 3   1st Qu.: 920000  
##  Malawi  :15   1997   : 3   Median :1100000  
##                1998   : 3   Mean   :1408889  
##                1999   : 3   3rd Qu.:2000000  
##                2000   : 3   Max.   :2600000  
##                (Other):27                    
##    population       
##  Min.   :9.964e+06  
##  1st Qu.:1.331e+07  
##  Median :6.995e+07  
##  Mean   :3.859e+08  
##  3rd Qu.:1.008e+09  
##  Max.   :1.190e+09  
## </code></pre>
<p><strong>Complete the code below to create a scatter plot of HIV incidence over time.</strong></p>
<p>Hint: You plot should look like this:</p>
<p><img src="https://drive.google.com/uc?export=view&id=1WPJkQBDW6jGnSTS3zRw3GGfNvEeDaLEr" width="480" /></p>
<pre class="r"><code>ggplot(data = {hiv_prevalence}, 
       {mapping} = aes({x} = year,
                       y = {total_cases})) {+}
  geom{_point}()</code></pre>
<hr />
<p>Make this plot and save it as an R object:
<img src="https://drive.google.com/uc?export=view&id=1hQNxU4hIFrNjhwk2DfKbsz7hpW6llXPZ" width="672" /></p>
<p>Put your plot object through the <code>ggplot_digest()</code> function and enter the resulting character string here: {3b276234cebfd6a5d17fde35457e039b}</p>
</div>
'@

$ws.Range("K2").Value = $k2
$ws.Range("K3").Value = $k3
$ws.Range("K4").Value = $k4

Write-Host "Updated K2, K3, K4"
